# Update column G ("K" = strikeouts) values per the regenerated save_data.
# The commit message indicates the save_data was regenerated to use K
# (strikeouts) instead of Strike# (pitch-level strike count), so the
# values in column G for most data rows are being recalculated/replaced.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$newValues = @{
    2 = 1
    3 = 0
    4 = 1
    5 = 0
    6 = 0
    7 = 0
    8 = 0
    9 = 1
    10 = 1
    11 = 1
    13 = 0
    14 = 0
    15 = 1
    16 = 0
    17 = 2
    18 = 1
    19 = 1
    20 = 1
    21 = 1
    23 = 0
    24 = 0
    25 = 2
    26 = 1
    27 = 1
    28 = 0
    29 = 0
    30 = 2
    31 = 1
    32 = 1
    33 = 0
    34 = 1
    35 = 1
    36 = 1
    37 = 1
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 0
    43 = 0
    44 = 0
    45 = 2
    46 = 2
    47 = 0
    48 = 2
    49 = 1
    50 = 2
    51 = 1
    52 = 1
    53 = 1
    54 = 1
    55 = 0
    56 = 0
    57 = 0
    58 = 2
    59 = 2
    60 = 0
    61 = 1
    62 = 1
    63 = 2
    64 = 1
    65 = 1
    66 = 0
    67 = 0
    68 = 0
    69 = 0
    70 = 2
    71 = 0
    72 = 2
    73 = 0
    74 = 1
    76 = 1
    77 = 0
    78 = 0
    79 = 1
    80 = 0
    81 = 1
    82 = 1
}

foreach ($row in $newValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $newValues[$row]
}
